$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 871235.25
$ws.Range("I15").Value = 871235.25
$ws.Range("K15").Value = 2613705.75
$ws.Range("M15").Value = -2613536.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2863.5
$ws.Range("I5").Value = 1700.3334
$ws.Range("J5").Value = 4026.6667
$ws.Range("K5").Value = 1700.3334
$ws.Range("L5").Value = 4026.6667
$ws.Range("M5").Value = -1588.3334
$ws.Range("N5").Value = -4250.6667
$ws.Range("H61").Value = 34483850
$ws.Range("I61").Value = 38462636
$ws.Range("J61").Value = 999.3333
$ws.Range("K61").Value = 38462636
$ws.Range("L61").Value = 999.3333
$ws.Range("M61").Value = -38462424
$ws.Range("N61").Value = -1423.3333
$ws.Range("H74").Value = 58829532
$ws.Range("I74").Value = 66672416
$ws.Range("K74").Value = 66672416
$ws.Range("M74").Value = -66671542
$ws.Range("H77").Value = 58829532
$ws.Range("I77").Value = 66672416
$ws.Range("K77").Value = 333362080
$ws.Range("M77").Value = -333357712
$ws.Range("H122").Value = 4258.4614
$ws.Range("I122").Value = 4067.25
$ws.Range("J122").Value = 4564.4
$ws.Range("K122").Value = 12201.75
$ws.Range("L122").Value = 13693.2
$ws.Range("M122").Value = -9751.75
$ws.Range("N122").Value = -18593.2
$ws.Range("H132").Value = 5885395
$ws.Range("I132").Value = 7145015.5
$ws.Range("K132").Value = 21435046.5
$ws.Range("M132").Value = -21432516.5
$ws.Range("H136").Value = 34483850
$ws.Range("I136").Value = 38462636
$ws.Range("J136").Value = 999.3333
$ws.Range("K136").Value = 115387908
$ws.Range("L136").Value = 2997.9999
$ws.Range("M136").Value = -115385358
$ws.Range("N136").Value = -8097.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2863.5
$ws.Range("I4").Value = 1700.3334
$ws.Range("J4").Value = 4026.6667
$ws.Range("K4").Value = 1700.3334
$ws.Range("L4").Value = 4026.6667
$ws.Range("M4").Value = -1585.3334
$ws.Range("N4").Value = -4256.6667
$ws.Range("H6").Value = 37952.668
$ws.Range("I6").Value = 30647
$ws.Range("J6").Value = 41605.5
$ws.Range("K6").Value = 30647
$ws.Range("L6").Value = 41605.5
$ws.Range("M6").Value = -30534
$ws.Range("N6").Value = -41831.5
$ws.Range("H22").Value = 2583.3076
$ws.Range("I22").Value = 3073.5
$ws.Range("J22").Value = 1799
$ws.Range("K22").Value = 3073.5
$ws.Range("L22").Value = 1799
$ws.Range("M22").Value = -2900.5
$ws.Range("N22").Value = -2145
$ws.Range("H95").Value = 16124.75
$ws.Range("J95").Value = 16124.75
$ws.Range("L95").Value = 16124.75
$ws.Range("N95").Value = -21616.75
$ws.Range("H109").Value = 66666
$ws.Range("J109").Value = 66666
$ws.Range("L109").Value = 66666
$ws.Range("N109").Value = -69440
$ws.Range("H111").Value = 55850
$ws.Range("J111").Value = 55850
$ws.Range("L111").Value = 55850
$ws.Range("N111").Value = -64030
$ws.Range("H112").Value = 58999
$ws.Range("J112").Value = 58999
$ws.Range("L112").Value = 58999
$ws.Range("N112").Value = -61953
$ws.Range("H117").Value = 31491
$ws.Range("J117").Value = 31491
$ws.Range("L117").Value = 31491
$ws.Range("N117").Value = -40669
$ws.Range("H118").Value = 183999.5
$ws.Range("J118").Value = 183999.5
$ws.Range("L118").Value = 183999.5
$ws.Range("N118").Value = -187313.5
$ws.Range("H134").Value = 15629613
$ws.Range("J134").Value = 2626
$ws.Range("L134").Value = 7878
$ws.Range("N134").Value = -12948

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8686.916999999999
$ws.Range("J22").Value = 623.25
$ws.Range("L22").Value = 623.25
$ws.Range("N22").Value = -1323.25
$ws.Range("H86").Value = 12887.25
$ws.Range("I86").Value = 9049.6
$ws.Range("J86").Value = 15628.429
$ws.Range("K86").Value = 9049.6
$ws.Range("L86").Value = 15628.429
$ws.Range("M86").Value = -7926.6
$ws.Range("N86").Value = -17874.429
$ws.Range("H89").Value = 12887.25
$ws.Range("I89").Value = 9049.6
$ws.Range("J89").Value = 15628.429
$ws.Range("K89").Value = 45248
$ws.Range("L89").Value = 78142.145
$ws.Range("M89").Value = -39632
$ws.Range("N89").Value = -89374.145
$ws.Range("H107").Value = 77642.766
$ws.Range("I107").Value = 585.8
$ws.Range("K107").Value = 585.8
$ws.Range("M107").Value = 1334.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 96.5
$ws.Range("I16").Value = 96.5
$ws.Range("K16").Value = 289.5
$ws.Range("M16").Value = -116.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1023624.75
$ws.Range("I11").Value = 1359999.6
$ws.Range("K11").Value = 1359999.6
$ws.Range("M11").Value = -1359860.6
$ws.Range("H13").Value = 144.4
$ws.Range("J13").Value = 144.4
$ws.Range("L13").Value = 144.4
$ws.Range("N13").Value = -422.4
$ws.Range("H102").Value = 3171.04
$ws.Range("I102").Value = 2445.8635
$ws.Range("K102").Value = 2445.8635
$ws.Range("M102").Value = -823.8634999999999
$ws.Range("H113").Value = 47182.918
$ws.Range("I113").Value = 53232.855
$ws.Range("K113").Value = 53232.855
$ws.Range("M113").Value = -51062.855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2410.5625
$ws.Range("I22").Value = 3053.625
$ws.Range("J22").Value = 1767.5
$ws.Range("K22").Value = 3053.625
$ws.Range("L22").Value = 1767.5
$ws.Range("M22").Value = -2758.625
$ws.Range("N22").Value = -2357.5
$ws.Range("H27").Value = 2410.5625
$ws.Range("I27").Value = 3053.625
$ws.Range("J27").Value = 1767.5
$ws.Range("K27").Value = 3053.625
$ws.Range("L27").Value = 1767.5
$ws.Range("M27").Value = -2946.625
$ws.Range("N27").Value = -1981.5
$ws.Range("H46").Value = 2124.8333
$ws.Range("I46").Value = 2124.8333
$ws.Range("K46").Value = 2124.8333
$ws.Range("M46").Value = -1936.8333
$ws.Range("H55").Value = 769.63635
$ws.Range("J55").Value = 1059
$ws.Range("L55").Value = 1059
$ws.Range("N55").Value = -1405
$ws.Range("H61").Value = 6325.375
$ws.Range("I61").Value = 6107.2856
$ws.Range("K61").Value = 6107.2856
$ws.Range("M61").Value = -5905.2856
$ws.Range("H92").Value = 21999
$ws.Range("J92").Value = 21999
$ws.Range("L92").Value = 21999
$ws.Range("N92").Value = -26991
$ws.Range("H113").Value = 6325.375
$ws.Range("I113").Value = 6107.2856
$ws.Range("K113").Value = 6107.2856
$ws.Range("M113").Value = -3937.2856
$ws.Range("H132").Value = 25266868
$ws.Range("I132").Value = 28238618
$ws.Range("K132").Value = 84715854
$ws.Range("M132").Value = -84713324

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4587.5
$ws.Range("I126").Value = 4587.5
$ws.Range("K126").Value = 13762.5
$ws.Range("M126").Value = -11292.5
$ws.Range("H136").Value = 14707229
$ws.Range("I136").Value = 14707229
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 44121687
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -44119137
$ws.Range("N136").ClearContents()
